$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 31   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/15/2024  Through  1/21/2024"

# --- Crime statistics table updates (rows 14-30) ---
$ws.Range("M14").Value = -100
$ws.Range("M14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("F15").Value = 2
$ws.Range("I15").Value = 1
$ws.Range("I15").NumberFormat = "#,##0"
$ws.Range("L15").Value = -66.666666666666
$ws.Range("N15").Value = -87.5
$ws.Range("C16").Value = 4
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "***.*"
$ws.Range("E16").NumberFormat = "General"
$ws.Range("F16").Value = 20
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = 53.846153846153
$ws.Range("I16").Value = 14
$ws.Range("J16").Value = 10
$ws.Range("K16").Value = 40
$ws.Range("L16").Value = 100
$ws.Range("M16").Value = -22.222222222222
$ws.Range("N16").Value = -90.728476821192
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -33.333333333333
$ws.Range("G17").Value = 35
$ws.Range("H17").Value = -28.571428571428
$ws.Range("I17").Value = 22
$ws.Range("J17").Value = 24
$ws.Range("K17").Value = -8.333333333333
$ws.Range("L17").Value = -12
$ws.Range("M17").Value = 4.761904761904
$ws.Range("N17").Value = -70.27027027027
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 20
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = 53.846153846153
$ws.Range("I18").Value = 15
$ws.Range("J18").Value = 8
$ws.Range("K18").Value = 87.5
$ws.Range("L18").Value = 50
$ws.Range("M18").Value = -28.571428571428
$ws.Range("N18").Value = -79.45205479452
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -42.857142857142
$ws.Range("F19").Value = 25
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = -16.666666666666
$ws.Range("I19").Value = 20
$ws.Range("J19").Value = 23
$ws.Range("K19").Value = -13.043478260869
$ws.Range("L19").Value = -20
$ws.Range("M19").Value = 53.846153846153
$ws.Range("N19").Value = -69.696969696969
$ws.Range("C20").Value = 4
$ws.Range("E20").Value = 300
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 71.428571428571
$ws.Range("I20").Value = 12
$ws.Range("J20").Value = 5
$ws.Range("K20").Value = 140
$ws.Range("L20").Value = 71.428571428571
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = -70.731707317073
$ws.Range("C21").Value = 21
$ws.Range("E21").Value = 5
$ws.Range("F21").Value = 104
$ws.Range("G21").Value = 98
$ws.Range("H21").Value = 6.122448979591
$ws.Range("I21").Value = 84
$ws.Range("J21").Value = 70
$ws.Range("K21").Value = 20
$ws.Range("L21").Value = 7.692307692307
$ws.Range("M21").Value = 5
$ws.Range("N21").Value = -79.661016949152
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("I22").Value = 1
$ws.Range("I22").NumberFormat = "#,##0"
$ws.Range("L22").Value = -50
$ws.Range("L22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M22").Value = -50
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 6
$ws.Range("E23").Value = -50
$ws.Range("F23").Value = 18
$ws.Range("G23").Value = 28
$ws.Range("H23").Value = -35.714285714285
$ws.Range("I23").Value = 15
$ws.Range("J23").Value = 19
$ws.Range("K23").Value = -21.052631578947
$ws.Range("L23").Value = -11.764705882352
$ws.Range("M23").Value = 650
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 29
$ws.Range("E24").Value = -31.03448275862
$ws.Range("F24").Value = 102
$ws.Range("G24").Value = 122
$ws.Range("H24").Value = -16.39344262295
$ws.Range("I24").Value = 71
$ws.Range("J24").Value = 88
$ws.Range("K24").Value = -19.318181818181
$ws.Range("L24").Value = -2.739726027397
$ws.Range("M24").Value = 29.090909090909
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 46
$ws.Range("G25").Value = 62
$ws.Range("H25").Value = -25.806451612903
$ws.Range("I25").Value = 34
$ws.Range("J25").Value = 45
$ws.Range("K25").Value = -24.444444444444
$ws.Range("L25").Value = -15
$ws.Range("M25").Value = -8.108108108108
$ws.Range("C26").Value = 1
$ws.Range("C26").NumberFormat = "#,##0"
$ws.Range("D26").Value = 1
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("E26").Value = 0
$ws.Range("E26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = -33.333333333333
$ws.Range("I26").Value = 1
$ws.Range("I26").NumberFormat = "#,##0"
$ws.Range("J26").Value = 3
$ws.Range("K26").Value = -66.666666666666
$ws.Range("L26").Value = -75
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("F27").Value = 5
$ws.Range("H27").Value = 66.666666666666
$ws.Range("J27").Value = 3
$ws.Range("K27").Value = 66.666666666666
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("C28").NumberFormat = "General"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "0"
$ws.Range("G28").NumberFormat = "General"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "***.*"
$ws.Range("H28").NumberFormat = "General"
$ws.Range("L28").Value = -50
$ws.Range("M28").Value = -50
$ws.Range("N28").Value = -95
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$ws.Range("C29").NumberFormat = "General"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "0"
$ws.Range("G29").NumberFormat = "General"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "***.*"
$ws.Range("H29").NumberFormat = "General"
$ws.Range("L29").Value = -50
$ws.Range("M29").Value = -50
$ws.Range("N29").Value = -93.75
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "***.*"
$ws.Range("E30").NumberFormat = "General"
$ws.Range("F30").NumberFormat = "@"
$ws.Range("F30").Value = "0"
$ws.Range("F30").NumberFormat = "General"
$ws.Range("H30").Value = -100
